$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the D (Price) column for this batch of edits so that
# numeric-looking strings (e.g. "649.71", "1.00", "0.165") are preserved as
# literal text instead of being auto-converted to numbers by COM automation.
$priceCol = $ws.Range("D2:D51")
$priceCol.NumberFormat = "@"

$ws.Range('D2').Value = '69.599.15'
$ws.Range('E2').Value = '  -0.28%  '
$ws.Range('D3').Value = '3.676.33'
$ws.Range('E3').Value = '  -0.85%  '
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('D5').Value = '649.71'
$ws.Range('E5').Value = '  -4.24%  '
$ws.Range('D6').Value = '160.83'
$ws.Range('E6').Value = '  -1.28%  '
$ws.Range('E7').Value = '  +0.17%  '
$ws.Range('D8').Value = '0.497'
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').Value = '0.145'
$ws.Range('E9').Value = '  -2.71%  '
$ws.Range('D10').Value = '7.16'
$ws.Range('E10').Value = '  -0.08%  '
$ws.Range('D11').Value = '0.442'
$ws.Range('E11').Value = '  -0.35%  '
$ws.Range('D12').Value = '0.0000231'
$ws.Range('E12').Value = '  -2.15%  '
$ws.Range('D13').Value = '4.304.92'
$ws.Range('E13').Value = '  -0.58%  '
$ws.Range('D14').Value = '32.62'
$ws.Range('E14').Value = '  -1.20%  '
$ws.Range('D15').Value = '3.699.30'
$ws.Range('E15').Value = '  -0.25%  '
$ws.Range('D16').Value = '69.690.11'
$ws.Range('E16').Value = '  -0.10%  '
$ws.Range('E17').Value = '  +0.38%  '
$ws.Range('D18').Value = '6.53'
$ws.Range('E18').Value = '  +0.20%  '
$ws.Range('D19').Value = '15.89'
$ws.Range('E19').Value = '  -1.53%  '
$ws.Range('D20').Value = '10.29'
$ws.Range('E20').Value = '  +4.29%  '
$ws.Range('D21').Value = '470.24'
$ws.Range('E21').Value = '  -0.77%  '
$ws.Range('D22').Value = '0.656'
$ws.Range('E22').Value = '  +0.14%  '
$ws.Range('D23').Value = '79.61'
$ws.Range('E23').Value = '  -1.12%  '
$ws.Range('D24').Value = '3.827.91'
$ws.Range('E24').Value = '  -0.68%  '
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').Value = '0.999'
$ws.Range('E25').Value = '  -0.12%  '
$ws.Range('B26').Value = 'PEPE'
$ws.Range('C26').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D26').Value = '0.0000126'
$ws.Range('E26').Value = '  -2.47%  '
$ws.Range('D27').Value = '11.18'
$ws.Range('E27').Value = '  +1.33%  '
$ws.Range('D28').Value = '8.81'
$ws.Range('E28').Value = '  -3.92%  '
$ws.Range('D29').Value = '2.65'
$ws.Range('E29').Value = '  -2.53%  '
$ws.Range('D30').Value = '1.70'
$ws.Range('E30').Value = '  -3.56%  '
$ws.Range('B31').Value = 'Binance-PegBSC-USD'
$ws.Range('C31').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D31').Value = '1.00'
$ws.Range('E31').Value = '  -0.13%  '
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').Value = '1.99'
$ws.Range('E32').Value = '  -2.19%  '
$ws.Range('D33').Value = '6.50'
$ws.Range('E33').Value = '  -2.10%  '
$ws.Range('B34').Value = 'Kaspa'
$ws.Range('C34').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D34').Value = '0.165'
$ws.Range('E34').Value = '  -0.48%  '
$ws.Range('D35').Value = '26.72'
$ws.Range('E35').Value = '  -1.10%  '
$ws.Range('D36').Value = '3.676.54'
$ws.Range('E36').Value = '  -0.49%  '
$ws.Range('D37').Value = '8.38'
$ws.Range('E37').Value = '  -2.54%  '
$ws.Range('E38').Value = '  -0.08%  '
$ws.Range('D39').Value = '5.87'
$ws.Range('E39').Value = '  -5.76%  '
$ws.Range('E40').Value = '  +0.21%  '
$ws.Range('D41').Value = '177.76'
$ws.Range('E41').Value = '  +6.06%  '
$ws.Range('B42').Value = 'Hedera'
$ws.Range('C42').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D42').Value = '0.0895'
$ws.Range('E42').Value = '  -1.42%  '
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').Value = '2.18'
$ws.Range('E43').Value = '  -2.72%  '
$ws.Range('D44').Value = '0.929'
$ws.Range('E44').Value = '  -1.90%  '
$ws.Range('B45').Value = 'OKB'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D45').Value = '46.75'
$ws.Range('E45').Value = '  -0.55%  '
$ws.Range('B46').Value = 'dogwifhat'
$ws.Range('C46').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D46').Value = '2.79'
$ws.Range('E46').Value = '  -0.29%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').Value = '29.02'
$ws.Range('E47').Value = '  +2.83%  '
$ws.Range('D48').Value = '0.000269'
$ws.Range('E48').Value = '  -5.01%  '
$ws.Range('D49').Value = '7.85'
$ws.Range('E49').Value = '  -1.10%  '
$ws.Range('D50').Value = '1.24'
$ws.Range('E50').Value = '  -4.62%  '
$ws.Range('E51').Value = '  -6.30%  '

# Restore the default (Normal) style on the price column so cells that were
# not touched, and the touched ones, keep the original unstyled appearance.
$priceCol.Style = "Normal"
